$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated snapshot of the crypto price table. Each entry is ($Cell, $Value).
# Column D ("Price") holds plain text in the source sheet (e.g. "26.837.07",
# "1.0000") so it gets a text number-format forced before the write (and
# reset back to the sheet's default style afterwards) to stop Excel from
# reinterpreting/rounding the literal as a number. Columns B/C/E are
# ordinary text and are written directly.
# (Rows 13/14 also swap which coin -- name (B) + link (C) -- occupies which
# rank; those four cells are included below like any other update.)
$updates = @(
    @{ Cell = 'D2'; Value = '26.837.07' },
    @{ Cell = 'E2'; Value = '  -1.13%  ' },
    @{ Cell = 'D3'; Value = '1.869.59' },
    @{ Cell = 'E3'; Value = '  -1.64%  ' },
    @{ Cell = 'E4'; Value = '  -0.28%  ' },
    @{ Cell = 'E5'; Value = '  -2.08%  ' },
    @{ Cell = 'D7'; Value = '0.5345' },
    @{ Cell = 'E7'; Value = '  +2.19%  ' },
    @{ Cell = 'D8'; Value = '0.3754' },
    @{ Cell = 'E8'; Value = '  -1.39%  ' },
    @{ Cell = 'D9'; Value = '0.07175' },
    @{ Cell = 'E9'; Value = '  -1.56%  ' },
    @{ Cell = 'D10'; Value = '21.59' },
    @{ Cell = 'E10'; Value = '  +1.04%  ' },
    @{ Cell = 'D11'; Value = '0.8859' },
    @{ Cell = 'E11'; Value = '  -1.89%  ' },
    @{ Cell = 'D12'; Value = '0.08133' },
    @{ Cell = 'E12'; Value = '  -0.43%  ' },
    @{ Cell = 'B13'; Value = 'Litecoin' },
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc' },
    @{ Cell = 'D13'; Value = '92.99' },
    @{ Cell = 'E13'; Value = '  -2.48%  ' },
    @{ Cell = 'B14'; Value = 'WrappedEther' },
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth' },
    @{ Cell = 'D14'; Value = '1.814.78' },
    @{ Cell = 'E14'; Value = '  -1.37%  ' },
    @{ Cell = 'D15'; Value = '5.264' },
    @{ Cell = 'E15'; Value = '  -1.71%  ' },
    @{ Cell = 'D16'; Value = '1.001' },
    @{ Cell = 'E16'; Value = '  -0.29%  ' },
    @{ Cell = 'D17'; Value = '14.73' },
    @{ Cell = 'E17'; Value = '  +0.32%  ' },
    @{ Cell = 'D18'; Value = '0.000008532' },
    @{ Cell = 'E18'; Value = '  -1.52%  ' },
    @{ Cell = 'D19'; Value = '1.000' },
    @{ Cell = 'E19'; Value = '  -0.20%  ' },
    @{ Cell = 'D20'; Value = '26.881.95' },
    @{ Cell = 'E20'; Value = '  -1.10%  ' },
    @{ Cell = 'D21'; Value = '4.969' },
    @{ Cell = 'E21'; Value = '  -3.11%  ' },
    @{ Cell = 'D22'; Value = '10.69' },
    @{ Cell = 'D23'; Value = '6.383' },
    @{ Cell = 'E23'; Value = '  -1.07%  ' },
    @{ Cell = 'D24'; Value = '147.11' },
    @{ Cell = 'E24'; Value = '  -1.40%  ' },
    @{ Cell = 'D25'; Value = '2.247' },
    @{ Cell = 'E25'; Value = '  -3.47%  ' },
    @{ Cell = 'D26'; Value = '1.731' },
    @{ Cell = 'E26'; Value = '  -0.53%  ' },
    @{ Cell = 'D27'; Value = '18.02' },
    @{ Cell = 'E27'; Value = '  -1.40%  ' },
    @{ Cell = 'E28'; Value = '  -1.15%  ' },
    @{ Cell = 'D29'; Value = '4.740' },
    @{ Cell = 'E29'; Value = '  -1.70%  ' },
    @{ Cell = 'D30'; Value = '4.574' },
    @{ Cell = 'E30'; Value = '  -6.59%  ' },
    @{ Cell = 'D31'; Value = '0.09132' },
    @{ Cell = 'E31'; Value = '  -0.98%  ' },
    @{ Cell = 'D32'; Value = '0.8003' },
    @{ Cell = 'E32'; Value = '  +1.01%  ' },
    @{ Cell = 'D33'; Value = '0.04990' },
    @{ Cell = 'E33'; Value = '  -1.14%  ' },
    @{ Cell = 'E34'; Value = '  +1.36%  ' },
    @{ Cell = 'D35'; Value = '1.167' },
    @{ Cell = 'E35'; Value = '  -4.45%  ' },
    @{ Cell = 'D36'; Value = '0.5921' },
    @{ Cell = 'E36'; Value = '  +3.58%  ' },
    @{ Cell = 'D37'; Value = '2.603' },
    @{ Cell = 'E37'; Value = '  -1.84%  ' },
    @{ Cell = 'D38'; Value = '3.152' },
    @{ Cell = 'E38'; Value = '  -6.24%  ' },
    @{ Cell = 'E39'; Value = '  -2.27%  ' },
    @{ Cell = 'E40'; Value = '  -1.27%  ' },
    @{ Cell = 'D41'; Value = '6.610' },
    @{ Cell = 'E41'; Value = '  +0.06%  ' },
    @{ Cell = 'D42'; Value = '8.858' },
    @{ Cell = 'E42'; Value = '  -2.41%  ' },
    @{ Cell = 'D43'; Value = '115.66' },
    @{ Cell = 'E43'; Value = '  -0.79%  ' },
    @{ Cell = 'D44'; Value = '0.5068' },
    @{ Cell = 'E44'; Value = '  +3.60%  ' },
    @{ Cell = 'D45'; Value = '0.1493' },
    @{ Cell = 'E45'; Value = '  -1.52%  ' },
    @{ Cell = 'D46'; Value = '1.0000' },
    @{ Cell = 'E46'; Value = '  -0.34%  ' },
    @{ Cell = 'D47'; Value = '9.919' },
    @{ Cell = 'E47'; Value = '  -2.16%  ' },
    @{ Cell = 'E48'; Value = '  -0.95%  ' },
    @{ Cell = 'D49'; Value = '37.63' },
    @{ Cell = 'E49'; Value = '  -2.25%  ' },
    @{ Cell = 'D50'; Value = '0.06027' },
    @{ Cell = 'E50'; Value = '  +1.21%  ' },
    @{ Cell = 'D51'; Value = '62.09' },
    @{ Cell = 'E51'; Value = '  -2.94%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Cell.StartsWith("D")) {
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
